$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) to the corrected/cleaned wording used when
# importing from file (extra spaces trimmed inconsistently upstream, and the
# currency unit added to the unit-price column).
$ws.Range("A1").Value = "Description "
$ws.Range("B1").Value = " Quantité"
$ws.Range("C1").Value = "Prix unitaire (€) "
$ws.Range("D1").Value = "TVA (%) "

# Widen the "Prix unitaire" column so the new "(€)" suffix fits.
$ws.Columns.Item(3).ColumnWidth = 15.42

# Missing/blank information coming from the import is no longer skipped: add
# a trailing row holding a blank placeholder cell instead of dropping it.
$ws.Cells.Item(9, 1).Value = " "

# Turn the data range into a real table (Table1) so the blank row above is
# recognised as part of the structured data when re-imported.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:D9"), $false, 1, "TableStyleMedium24")
$tbl.TableStyle = "TableStyleMedium24"
